$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that is refreshed for
# every data row whenever the sheet is regenerated. Bump it from
# 2023-09-10 (45179) to 2023-09-11 (45180) for every row (C2:C463).
$ws.Range("C2:C463").Value = 45180
